$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")

$a2 = $aboutSheet.Range("A2").Value()
$aboutSheet.Range("A2").Value = $a2.Replace($oldVersion, $newVersion)

$a6 = $aboutSheet.Range("A6").Value()
$aboutSheet.Range("A6").Value = $a6.Replace($oldVersion, $newVersion)

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 8; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # Column S = 19
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldVersion, $newVersion)
}
